$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# --- Update the "Date" value (row 8, column B) to the new timestamp ---
$ws.Range("B8").Value = "2025-11-04T10:04:56+00:00"

# --- Set the "Experimental" value (row 7, column B) to the literal text "false" ---
# A plain Value/Formula assignment of the bare word "false" gets auto-coerced to a
# Boolean by Excel (same as typing FALSE into a General-formatted cell). To store it
# as literal text (matching the target workbook) we compute it with a text formula
# in a scratch cell, copy the evaluated result, and paste-special as values into the
# target cell - this yields a genuine text cell instead of a boolean.
$scratch = $ws.Range("D1")
$scratch.Formula = "=""false"""
$scratch.Copy()
$ws.Range("B7").PasteSpecial(-4163)  # xlPasteValues
$scratch.ClearContents()
